$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$theme = $m.Theme
Write-Host $theme
Write-Host $theme.ThemeVariants.Count
$d = $p.Designs
Write-Host $d.Count
for ($i=1; $i -le $d.Count; $i++) {
    $item = $d.Item($i)
    Write-Host "Design $i : $($item.Name)"
}
